# Auto-generated Excel COM-interop script to apply Typhon_Profits value refresh
# Updates cached numeric values in columns H-N across multiple worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 849.6111
$ws.Cells.Item(28, 9).Value = 397.22223
$ws.Cells.Item(28, 10).Value = 1302
$ws.Cells.Item(28, 11).Value = 397.22223
$ws.Cells.Item(28, 12).Value = 1302
$ws.Cells.Item(28, 13).Value = 87.77776999999998
$ws.Cells.Item(28, 14).Value = -2272
$ws.Cells.Item(74, 8).Value = 9619423
$ws.Cells.Item(74, 9).Value = 3000
$ws.Cells.Item(74, 10).Value = 13893389
$ws.Cells.Item(74, 11).Value = 3000
$ws.Cells.Item(74, 12).Value = 13893389
$ws.Cells.Item(74, 13).Value = -2064
$ws.Cells.Item(74, 14).Value = -13895261
$ws.Cells.Item(77, 8).Value = 9619423
$ws.Cells.Item(77, 9).Value = 3000
$ws.Cells.Item(77, 10).Value = 13893389
$ws.Cells.Item(77, 11).Value = 15000
$ws.Cells.Item(77, 12).Value = 69466945
$ws.Cells.Item(77, 13).Value = -10320
$ws.Cells.Item(77, 14).Value = -69476305
$ws.Cells.Item(86, 8).Value = 5430.087
$ws.Cells.Item(86, 9).Value = 1103.125
$ws.Cells.Item(86, 11).Value = 1103.125
$ws.Cells.Item(86, 13).Value = 19.875
$ws.Cells.Item(89, 8).Value = 5430.087
$ws.Cells.Item(89, 9).Value = 1103.125
$ws.Cells.Item(89, 11).Value = 5515.625
$ws.Cells.Item(89, 13).Value = 100.375
$ws.Cells.Item(108, 8).Value = 35999
$ws.Cells.Item(108, 10).Value = 35999
$ws.Cells.Item(108, 12).Value = 35999
$ws.Cells.Item(108, 14).Value = -43679
$ws.Cells.Item(113, 8).Value = 4146.2
$ws.Cells.Item(113, 9).Value = 3126.6365
$ws.Cells.Item(113, 11).Value = 3126.6365
$ws.Cells.Item(113, 13).Value = 127.3634999999999
$ws.Cells.Item(116, 8).Value = 6455.091
$ws.Cells.Item(129, 8).Value = 313767.7
$ws.Cells.Item(129, 10).Value = 334650.53
$ws.Cells.Item(129, 12).Value = 1003951.59
$ws.Cells.Item(129, 14).Value = -1013951.59
$ws.Cells.Item(135, 8).Value = 16134277
$ws.Cells.Item(135, 9).Value = 641.7727
$ws.Cells.Item(135, 10).Value = 55572052
$ws.Cells.Item(135, 11).Value = 5775.954299999999
$ws.Cells.Item(135, 12).Value = 500148468
$ws.Cells.Item(135, 13).Value = -3240.954299999999
$ws.Cells.Item(135, 14).Value = -500153538
$ws.Cells.Item(141, 8).Value = 2268.2222
$ws.Cells.Item(141, 9).Value = 1910.0834
$ws.Cells.Item(141, 10).Value = 5133.3335
$ws.Cells.Item(141, 11).Value = 5730.2502
$ws.Cells.Item(141, 12).Value = 15400.0005
$ws.Cells.Item(141, 13).Value = -550.2502000000004
$ws.Cells.Item(141, 14).Value = -25760.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 1006
$ws.Cells.Item(16, 9).Value = 1006
$ws.Cells.Item(16, 11).Value = 1006
$ws.Cells.Item(16, 13).Value = -719
$ws.Cells.Item(32, 8).Value = 5922.339
$ws.Cells.Item(32, 9).Value = 4161.14
$ws.Cells.Item(32, 10).Value = 26000
$ws.Cells.Item(32, 11).Value = 4161.14
$ws.Cells.Item(32, 12).Value = 26000
$ws.Cells.Item(32, 13).Value = -3874.14
$ws.Cells.Item(32, 14).Value = -26574
$ws.Cells.Item(88, 8).Value = 101585.1
$ws.Cells.Item(88, 9).Value = 1530.75
$ws.Cells.Item(88, 10).Value = 168288
$ws.Cells.Item(88, 11).Value = 1530.75
$ws.Cells.Item(88, 12).Value = 168288
$ws.Cells.Item(88, 13).Value = -1124.75
$ws.Cells.Item(88, 14).Value = -169100
$ws.Cells.Item(91, 8).Value = 101585.1
$ws.Cells.Item(91, 9).Value = 1530.75
$ws.Cells.Item(91, 10).Value = 168288
$ws.Cells.Item(91, 11).Value = 1530.75
$ws.Cells.Item(91, 12).Value = 168288
$ws.Cells.Item(91, 13).Value = -126.75
$ws.Cells.Item(91, 14).Value = -171096
$ws.Cells.Item(97, 8).Value = 125001050
$ws.Cells.Item(97, 9).Value = 1101.6666
$ws.Cells.Item(97, 11).Value = 1101.6666
$ws.Cells.Item(97, 13).Value = -605.6666
$ws.Cells.Item(132, 8).Value = 21620.04
$ws.Cells.Item(132, 9).Value = 2158.3684
$ws.Cells.Item(132, 11).Value = 6475.1052
$ws.Cells.Item(132, 13).Value = -3945.1052
$ws.Cells.Item(133, 8).Value = 50000
$ws.Cells.Item(133, 10).Value = 50000
$ws.Cells.Item(133, 12).Value = 50000
$ws.Cells.Item(133, 14).Value = -55060

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2067.6
$ws.Cells.Item(99, 9).Value = 1795.8
$ws.Cells.Item(99, 11).Value = 1795.8
$ws.Cells.Item(99, 13).Value = -297.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3627.5278
$ws.Cells.Item(31, 9).Value = 3374.9092
$ws.Cells.Item(31, 10).Value = 3738.68
$ws.Cells.Item(31, 11).Value = 3374.9092
$ws.Cells.Item(31, 12).Value = 3738.68
$ws.Cells.Item(31, 13).Value = -3079.9092
$ws.Cells.Item(31, 14).Value = -4328.68
$ws.Cells.Item(34, 8).Value = 3627.5278
$ws.Cells.Item(34, 9).Value = 3374.9092
$ws.Cells.Item(34, 10).Value = 3738.68
$ws.Cells.Item(34, 11).Value = 3374.9092
$ws.Cells.Item(34, 12).Value = 3738.68
$ws.Cells.Item(34, 13).Value = -3172.9092
$ws.Cells.Item(34, 14).Value = -4142.68
$ws.Cells.Item(62, 8).Value = 52635590
$ws.Cells.Item(62, 10).Value = 4716
$ws.Cells.Item(62, 12).Value = 4716
$ws.Cells.Item(62, 14).Value = -5964
$ws.Cells.Item(65, 8).Value = 52635590
$ws.Cells.Item(65, 10).Value = 4716
$ws.Cells.Item(65, 12).Value = 23580
$ws.Cells.Item(65, 14).Value = -29820
$ws.Cells.Item(134, 8).Value = 1447.0588
$ws.Cells.Item(134, 9).Value = 1412.5
$ws.Cells.Item(134, 11).Value = 4237.5
$ws.Cells.Item(134, 13).Value = -1702.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 79.36364
$ws.Cells.Item(12, 9).Value = 6.6666665
$ws.Cells.Item(12, 10).Value = 106.625
$ws.Cells.Item(12, 11).Value = 19.9999995
$ws.Cells.Item(12, 12).Value = 319.875
$ws.Cells.Item(12, 13).Value = 153.0000005
$ws.Cells.Item(12, 14).Value = -665.875
$ws.Cells.Item(26, 8).Value = 600.1667
$ws.Cells.Item(26, 9).Value = 520.2
$ws.Cells.Item(26, 10).Value = 1000
$ws.Cells.Item(26, 11).Value = 1560.6
$ws.Cells.Item(26, 12).Value = 3000
$ws.Cells.Item(26, 13).Value = -1272.6
$ws.Cells.Item(26, 14).Value = -3576
$ws.Cells.Item(128, 8).Value = 190000
$ws.Cells.Item(128, 9).Value = 190000
$ws.Cells.Item(128, 11).Value = 570000
$ws.Cells.Item(128, 13).Value = -565020
$ws.Cells.Item(131, 8).Value = 725.11
$ws.Cells.Item(131, 9).Value = 515
$ws.Cells.Item(131, 10).Value = 729.39795
$ws.Cells.Item(131, 11).Value = 1545
$ws.Cells.Item(131, 12).Value = 2188.19385
$ws.Cells.Item(131, 13).Value = 3495
$ws.Cells.Item(131, 14).Value = -12268.19385
$ws.Cells.Item(140, 8).Value = 1903.75
$ws.Cells.Item(140, 9).Value = 1405.4546
$ws.Cells.Item(140, 10).Value = 3000
$ws.Cells.Item(140, 11).Value = 4216.3638
$ws.Cells.Item(140, 12).Value = 9000
$ws.Cells.Item(140, 13).Value = 963.6361999999999
$ws.Cells.Item(140, 14).Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 58.1875
$ws.Cells.Item(2, 9).Value = 55.583332
$ws.Cells.Item(2, 10).Value = 66
$ws.Cells.Item(2, 11).Value = 55.583332
$ws.Cells.Item(2, 12).Value = 66
$ws.Cells.Item(2, 13).Value = 57.416668
$ws.Cells.Item(2, 14).Value = -292
$ws.Cells.Item(122, 8).Value = 78432720
$ws.Cells.Item(122, 10).Value = 166668110
$ws.Cells.Item(122, 12).Value = 500004330
$ws.Cells.Item(122, 14).Value = -500009230
$ws.Cells.Item(132, 8).Value = 28855.684
$ws.Cells.Item(132, 9).Value = 2091.3635
$ws.Cells.Item(132, 11).Value = 6274.0905
$ws.Cells.Item(132, 13).Value = -3744.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 546
$ws.Cells.Item(4, 9).Value = 546
$ws.Cells.Item(4, 11).Value = 546
$ws.Cells.Item(4, 13).Value = -433
$ws.Cells.Item(24, 8).Value = 5666.6665
$ws.Cells.Item(24, 10).Value = 5666.6665
$ws.Cells.Item(24, 12).Value = 5666.6665
$ws.Cells.Item(24, 14).Value = -6352.6665
$ws.Cells.Item(28, 8).Value = 546
$ws.Cells.Item(28, 9).Value = 546
$ws.Cells.Item(28, 11).Value = 546
$ws.Cells.Item(28, 13).Value = -314
$ws.Cells.Item(37, 8).Value = 546
$ws.Cells.Item(37, 9).Value = 546
$ws.Cells.Item(37, 11).Value = 546
$ws.Cells.Item(37, 13).Value = -439
$ws.Cells.Item(104, 8).Value = 17768.4
$ws.Cells.Item(104, 10).Value = 17768.4
$ws.Cells.Item(104, 12).Value = 17768.4
$ws.Cells.Item(104, 14).Value = -24756.4
$ws.Cells.Item(122, 8).Value = 894255.0600000001
$ws.Cells.Item(122, 10).Value = 3220.6155
$ws.Cells.Item(122, 12).Value = 9661.8465
$ws.Cells.Item(122, 14).Value = -14561.8465
$ws.Cells.Item(132, 8).Value = 1941.2307
$ws.Cells.Item(132, 9).Value = 1390.3636
$ws.Cells.Item(132, 11).Value = 4171.0908
$ws.Cells.Item(132, 13).Value = -1641.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 200000
$ws.Cells.Item(14, 9).Value = 200000
$ws.Cells.Item(14, 11).Value = 200000
$ws.Cells.Item(14, 13).Value = -199832
$ws.Cells.Item(62, 8).Value = 2977.4443
$ws.Cells.Item(62, 9).Value = 2500
$ws.Cells.Item(62, 10).Value = 3113.8572
$ws.Cells.Item(62, 11).Value = 2500
$ws.Cells.Item(62, 12).Value = 3113.8572
$ws.Cells.Item(62, 13).Value = -1876
$ws.Cells.Item(62, 14).Value = -4361.8572
$ws.Cells.Item(65, 8).Value = 2977.4443
$ws.Cells.Item(65, 9).Value = 2500
$ws.Cells.Item(65, 10).Value = 3113.8572
$ws.Cells.Item(65, 11).Value = 12500
$ws.Cells.Item(65, 12).Value = 15569.286
$ws.Cells.Item(65, 13).Value = -9380
$ws.Cells.Item(65, 14).Value = -21809.286
$ws.Cells.Item(74, 8).Value = 31265.4
$ws.Cells.Item(74, 10).Value = 31265.4
$ws.Cells.Item(74, 12).Value = 31265.4
$ws.Cells.Item(74, 14).Value = -33137.4
$ws.Cells.Item(77, 8).Value = 31265.4
$ws.Cells.Item(77, 10).Value = 31265.4
$ws.Cells.Item(77, 12).Value = 93796.20000000001
$ws.Cells.Item(77, 14).Value = -103156.2
$ws.Cells.Item(122, 8).Value = 846.97675
$ws.Cells.Item(122, 10).Value = 957.93335
$ws.Cells.Item(122, 12).Value = 2873.80005
$ws.Cells.Item(122, 14).Value = -7773.80005
